# Pixell test plan - ChequingAccount class unit test plan
# Fills in developer name and the test-case detail columns (Method Inputs /
# Condition being Tested / Expected Result) for test cases 1-8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Developer name
$ws.Range("C3").Value = "Raven Manalastas"

# Test case 1 - __init__ - Attributes are set to input values (ensure to
# test for superclass and subclass attributes)
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = "Account Number = 710`nClient Number = 7910`nBlance = 2500.00`nDate Created = July 10, 2024`noverdraft_limit = 50`noverdraft_rate = .10"
$ws.Range("G7").Value = "Attributes are set to input values."

# Test case 2 - __init__ - overdraft limit has invalid type.
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = 'overdraft_limit = "Invalid limit"'
$ws.Range("G8").Value = "overdraft_limit attribute is set to -100"

# Test case 3 - __init__ - overdraft rate has invalid type.
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = 'overdraft_rate = "Invalid rate"'
$ws.Range("G9").Value = "overdraft_rate attribute is set to 0.05"

# Test case 4 - __init__ - date created has invalid type
$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = 'date_created: = "Invalid date"'
$ws.Range("G10").Value = "date_created attribute uses today()method"

# Test case 5 - __str__ - balance greater than overdraft limit
$ws.Range("E11").Value = "None"
$ws.Range("F11").Value = "balance = 2500"
$ws.Range("G11").Value = "calculated_service_charge = 0.50"

# Test case 6 - __str__ - balance less than overdraft limit
$ws.Range("E12").Value = "None"
$ws.Range("F12").Value = "balance = -300"
$ws.Range("G12").Value = "calculated_service_charge = 10.50"

# Test case 7 - __str__ - balance equal to overdraft limit
$ws.Range("E13").Value = "None"
$ws.Range("F13").Value = "balance = 50"
$ws.Range("G13").Value = "calculated_service_charge = 0.50"

# Test case 8 - get_service_charges - appropriate value returned based on
# attribute values.
$ws.Range("E14").Value = "None"
$ws.Range("F14").Value = "Account Number = 710`nClient Number = 7910`nBlance = 2500.00`nDate Created = July 10, 2024`noverdraft_limit = 50`noverdraft_rate = .10"
$ws.Range("G14").Value = "call the str method of the superclass BankAccount and the then concatenate with a formatted f string.`nf`"Account Number: {self.__account_number} Balance: {self.__balance} `n`"`nf`"Overdraft Limit: `${self.__overdraft_limit:,.2f} Overdraft Rate: {self.__overdraft_rate*100:.2f} % Account Type: Chequing`""

# Window/selection state left by the editing session
$ws.Range("H13").Select()
$excel.ActiveWindow.Zoom = 84
